$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.792.79"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "1.887.50"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "0.7924"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").Value = "241.98"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "0.3172"
$ws.Range("E8").Value = "  +2.01%  "
$ws.Range("D9").Value = "25.46"
$ws.Range("E9").Value = "  -3.59%  "
$ws.Range("D10").Value = "0.07040"
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("D11").Value = "0.08050"
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").Value = "0.7660"
$ws.Range("E12").Value = "  +3.73%  "
$ws.Range("D13").Value = "1.895.90"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").Value = "5.292"
$ws.Range("E14").Value = "  +2.50%  "
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "29.813.80"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "13.81"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").Value = "5.927"
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("D19").Value = "243.25"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").Value = "0.000007717"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "2.134.14"
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("D23").Value = "8.076"
$ws.Range("E23").Value = "  +17.18%  "
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "0.1613"
$ws.Range("E25").Value = "  +11.04%  "
$ws.Range("D26").Value = "9.290"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("D27").Value = "163.98"
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("D28").Value = "18.65"
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("D29").Value = "2.053"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "1.368"
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("D31").Value = "1.535"
$ws.Range("E31").Value = "  +1.57%  "
$ws.Range("D32").Value = "4.433"
$ws.Range("E32").Value = "  +3.68%  "
$ws.Range("D33").Value = "0.05630"
$ws.Range("E33").Value = "  +2.16%  "
$ws.Range("D34").Value = "4.090"
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("D35").Value = "1.263"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").Value = "0.9981"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "2.716"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").Value = "0.01924"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").Value = "0.4416"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("D42").Value = "72.11"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "5.834"
$ws.Range("E43").Value = "  -2.21%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "0.8404"
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("D46").Value = "1.027.31"
$ws.Range("E46").Value = "  +4.99%  "
$ws.Range("D47").Value = "1.869"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").Value = "101.84"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("D49").Value = "9.922"
$ws.Range("E49").Value = "  +2.85%  "
$ws.Range("D50").Value = "7.444"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").Value = "2.049.21"
$ws.Range("E51").Value = "  -0.90%  "
